$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump Version value, then remove the Jurisdiction/Chile row ---
$meta = $wb.Worksheets.Item("Metadata")

# Update the Version property value (row 3, column B) from 0.4.0 to 0.7.0
$meta.Range("B3").Value = "0.7.0"

# Remove the entire "Jurisdiction" / "Chile" row (row 11); rows below shift up.
$meta.Rows(11).Delete()
